$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to store the value as literal text so that
    # numeric-looking strings (e.g. "1.00", "0.671") are not
    # auto-converted to numbers and lose formatting/precision.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$sub3 = [char]0x2083

$ws.Range("D2").Value = "70.587.68"
$ws.Range("E2").Value = "  -3.09%  "
$ws.Range("D3").Value = "3.853.72"
$ws.Range("E3").Value = "  -3.21%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.13%  "
Set-TextValue $ws.Range("D5") "591.49"
$ws.Range("E5").Value = "  -0.05%  "
Set-TextValue $ws.Range("D6") "166.21"
$ws.Range("E6").Value = "  +4.23%  "
Set-TextValue $ws.Range("D7") "0.671"
$ws.Range("E7").Value = "  -2.09%  "
$ws.Range("E8").Value = "  +0.22%  "
Set-TextValue $ws.Range("D9") "0.747"
$ws.Range("E9").Value = "  -0.45%  "
Set-TextValue $ws.Range("D10") "0.175"
$ws.Range("E10").Value = "  +4.10%  "
Set-TextValue $ws.Range("D11") "53.04"
$ws.Range("E11").Value = "  -1.64%  "
Set-TextValue $ws.Range("D12") "0.0000319"
$ws.Range("E12").Value = "  +0.13%  "
Set-TextValue $ws.Range("D13") "11.36"
$ws.Range("E13").Value = "  +3.99%  "
$ws.Range("D14").Value = "4.475.15"
$ws.Range("E14").Value = "  -2.91%  "
Set-TextValue $ws.Range("D15") "21.11"
$ws.Range("E15").Value = "  +3.58%  "
$ws.Range("D16").Value = "3.876.08"
$ws.Range("E16").Value = "  -2.63%  "
Set-TextValue $ws.Range("D17") "13.76"
$ws.Range("E17").Value = "  -2.26%  "
$ws.Range("E18").Value = "  -5.66%  "
$ws.Range("E19").Value = "  -2.19%  "
$ws.Range("D20").Value = "70.561.16"
$ws.Range("E20").Value = "  -2.79%  "
Set-TextValue $ws.Range("D21") "437.18"
$ws.Range("E21").Value = "  +0.34%  "
Set-TextValue $ws.Range("D22") "4.71"
$ws.Range("E22").Value = "  -1.25%  "
Set-TextValue $ws.Range("D23") "93.80"
$ws.Range("E23").Value = "  -2.41%  "
Set-TextValue $ws.Range("D24") "3.25"
$ws.Range("E24").Value = "  -5.01%  "
Set-TextValue $ws.Range("D25") "13.82"
$ws.Range("E25").Value = "  -3.45%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D26") "4.04"
$ws.Range("E26").Value = "  -9.76%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D27") "11.15"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D28") "5.93"
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D29") "10.39"
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D30") "35.05"
$ws.Range("E30").Value = "  -3.62%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D31") "8.14"
$ws.Range("E31").Value = "  +4.07%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D32") "13.47"
$ws.Range("E32").Value = "  -1.83%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D33") "48.15"
$ws.Range("E33").Value = "  -0.20%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D34") "0.125"
$ws.Range("E34").Value = "  -4.82%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D35") "69.64"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0${sub3}0975"
$ws.Range("E36").Value = "  +10.97%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D37") "632.96"
$ws.Range("E37").Value = "  -6.76%  "
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D38") "0.421"
$ws.Range("E38").Value = "  -3.23%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D39") "0.145"
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D40") "1.00"
$ws.Range("E40").Value = "  +0.02%  "
Set-TextValue $ws.Range("D41") "1.00"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D42") "3.29"
$ws.Range("E42").Value = "  -3.06%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D43") "3.28"
$ws.Range("E43").Value = "  +25.09%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D44") "0.0469"
$ws.Range("E44").Value = "  -3.69%  "
$ws.Range("B45").Value = "THORChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D45") "10.02"
$ws.Range("E45").Value = "  -7.63%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D46") "2.68"
$ws.Range("E46").Value = "  +1.40%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D47") "0.143"
$ws.Range("E47").Value = "  -4.32%  "
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D48") "2.84"
$ws.Range("E48").Value = "  -15.04%  "
$ws.Range("D49").Value = "2.830.21"
$ws.Range("E49").Value = "  +1.38%  "
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Range("D50") "3.22"
$ws.Range("E50").Value = "  -5.43%  "
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextValue $ws.Range("D51") "0.000271"
$ws.Range("E51").Value = "  +0.97%  "
